$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 26
$ws.Range("B2").Value = 50
$ws.Range("B3").Value = 108
$ws.Range("B4").Value = 145
$ws.Range("B5").Value = 190
$ws.Range("B6").Value = 234
